$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B3").Value = "1.1.0"
$metadata.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$metadata.Range("B10").Value = "No display for ContactDetail"

# --- Include from FSIII sheet updates ---
$include1 = $wb.Worksheets.Item("Include from FSIII")
$include1.Range("C2").Value = "A"

# --- Remove the extra "Include from FSIII 2" sheet (revert) ---
$include2 = $wb.Worksheets.Item("Include from FSIII 2")
$include2.Delete()
